# GHGs through 8 Dec 2021; corrected megasheet.
# Append the new Summer/Fall 2021 headspace-prep rows (36-43) to Sheet1 and
# apply the "final row" box-border formatting to the last row of data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- new data rows 36-41 (same plain style as the rows above them) ----
$ws.Range("A36").Value = 44488
$ws.Range("B36").Value = 21
$ws.Range("C36").Value = 30.234000000000002

$ws.Range("A37").Value = 44495
$ws.Range("B37").Value = 19.399999999999999
$ws.Range("C37").Value = 29.742000000000001

$ws.Range("A38").Value = 44496
$ws.Range("B38").Value = 20.399999999999999
$ws.Range("C38").Value = 29.939

$ws.Range("A39").Value = 44504
$ws.Range("B39").Value = 20.9
$ws.Range("C39").Value = 30.305

$ws.Range("A40").Value = 44509
$ws.Range("B40").Value = 20.8
$ws.Range("C40").Value = 30.274999999999999

$ws.Range("A41").Value = 44510
$ws.Range("B41").Value = 20.2
$ws.Range("C41").Value = 30.169

# ---- row 42: last row of the "old" formatting block ----
$ws.Range("A42").Value = 44523
$ws.Range("B42").Value = 20.8
$ws.Range("C42").Value = 30.228000000000002

# ---- row 43: newest reading, boxed with a medium light-grey border and
#      right-aligned / wrapped text, marking it as the sheet's last row ----
$ws.Range("A43").Value = 44537
$ws.Range("B43").Value = 20.3
$ws.Range("C43").Value = 30.210999999999999

$lastRow = $ws.Range("A43:C43")
$lastRow.HorizontalAlignment = -4152
$lastRow.WrapText = $true
$lastRow.Borders.Weight = -4138
$lastRow.Borders.Color = 13421772

# Date formatting for the new date cells (matches the rest of column A)
$ws.Range("A36:A43").NumberFormat = "d-mmm-yy"

# ---- view state: selection on the newly added last row ----
$ws.Range("A43:C43").Select()
$excel.ActiveWindow.ScrollRow = 15

# ---- page setup ----
$ws.PageSetup.Orientation = 1

Write-Host "Applied GHG update through 8 Dec 2021."
